# Apply weekly fruit/vegetable price update:
# Insert a new record row at row 269 (shifting existing rows 269-287 down to 270-288)
# and populate it with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 269, pushing existing data down
$ws.Rows.Item(269).Insert()

# Populate the new row 269 with the new record
$ws.Cells.Item(269, 1).Value = 9
$ws.Cells.Item(269, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(269, 3).Value = "Metropolitana"
$ws.Cells.Item(269, 4).Value = 44714
$ws.Cells.Item(269, 5).Value = 13
$ws.Cells.Item(269, 6).Value = 100112043
$ws.Cells.Item(269, 7).Value = "Pepino ensalada"
$ws.Cells.Item(269, 8).Value = "Sin especificar"
$ws.Cells.Item(269, 9).Value = "Primera"
$ws.Cells.Item(269, 10).Value = 61
$ws.Cells.Item(269, 11).Value = 16000
$ws.Cells.Item(269, 12).Value = 18000
$ws.Cells.Item(269, 13).Value = 17016
$ws.Cells.Item(269, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(269, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(269, 16).Value = 284
$ws.Cells.Item(269, 17).Value = 60
$ws.Cells.Item(269, 18).Value = "Hortaliza"
